# Citrix.WEMSDK.xlsx — "Get-WEMCondition & New-WEMCondition completed"
#
# Adds the new Filter "Condition" and "Rule" command rows to the
# "Commands and Aliasses" sheet (rows 61-68). Get-WEMCondition and
# New-WEMCondition are marked "Ready for testing" (completed); the
# remaining Condition/Rule cmdlets are still "In Development".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands and Aliasses")

# Row 61: Get-WEMCondition (completed -> Ready for testing)
$ws.Cells.Item(61, 1).Value = "Ready for testing"
$ws.Cells.Item(61, 2).Value = "Get-WEMCondition"
$ws.Cells.Item(61, 4).Value = "Returns one or more Filter Condition objects from the WEM Database."

# Row 62: Set-WEMCondition (still in development)
$ws.Cells.Item(62, 1).Value = "In Development"
$ws.Cells.Item(62, 2).Value = "Set-WEMCondition"
$ws.Cells.Item(62, 4).Value = "Updates a Filter Condition object from the WEM Database."

# Row 63: New-WEMCondition (completed -> Ready for testing)
$ws.Cells.Item(63, 1).Value = "Ready for testing"
$ws.Cells.Item(63, 2).Value = "New-WEMCondition"
$ws.Cells.Item(63, 4).Value = "Creates a Filter Condition object from the WEM Database."

# Row 64: Remove-WEMCondition (still in development)
$ws.Cells.Item(64, 1).Value = "In Development"
$ws.Cells.Item(64, 2).Value = "Remove-WEMCondition"
$ws.Cells.Item(64, 4).Value = "Removes a Filter Condition object from the WEM Database."

# Row 65: Get-WEMRule (still in development)
$ws.Cells.Item(65, 1).Value = "In Development"
$ws.Cells.Item(65, 2).Value = "Get-WEMRule"
$ws.Cells.Item(65, 4).Value = "Returns one or more Filter Rule objects from the WEM Database."

# Row 66: Set-WEMRule (still in development)
$ws.Cells.Item(66, 1).Value = "In Development"
$ws.Cells.Item(66, 2).Value = "Set-WEMRule"
$ws.Cells.Item(66, 4).Value = "Updates a Filter Rule object from the WEM Database."

# Row 67: New-WEMRule (still in development)
$ws.Cells.Item(67, 1).Value = "In Development"
$ws.Cells.Item(67, 2).Value = "New-WEMRule"
$ws.Cells.Item(67, 4).Value = "Creates a Filter Rule object from the WEM Database."

# Row 68: Remove-WEMRule (still in development)
$ws.Cells.Item(68, 1).Value = "In Development"
$ws.Cells.Item(68, 2).Value = "Remove-WEMRule"
$ws.Cells.Item(68, 4).Value = "Removes a Filter Rule object from the WEM Database."

# Reflect the author's last on-screen position: scrolled down with
# cell B63 (New-WEMCondition) selected.
$ws.Activate()
$ws.Range("B63").Select()

Write-Output "Added Filter Condition/Rule command rows (61-68)."
